$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row cells from human-readable labels to camelCase field
# names (to match normalized naming used by the naive bayes classifier /
# front-end prediction page), and shorten the last header's label.
$ws.Range("A1").Value = "jenisKelamin"
$ws.Range("B1").Value = "organisasi"
$ws.Range("C1").Value = "ekstrakurikuler"
$ws.Range("D1").Value = "sertifikasiProfesi"
$ws.Range("E1").Value = "nilaiAkhir"
$ws.Range("F1").Value = "tempatMagang"
$ws.Range("G1").Value = "tempatKerja"
$ws.Range("H1").Value = "Durasi Mendapat Kerja"
